$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetimes
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-22 06:30:08"
$wsZhCn.Range("G2").Value = "2016-02-22 06:31:06"

# de-de sheet: update handoff/handback datetimes
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-22 06:30:22"
$wsDeDe.Range("G2").Value = "2016-02-22 06:31:33"
